$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Characters(21, 2).Text = "34"
$ws.Range("C9").Characters(27, 9).Text = "8/21/2023"
$ws.Range("C9").Characters(47, 9).Text = "8/27/2023"

# --- Crime-stat table updates (rows 15-29) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("L15").Value = -75
$ws.Range("C16").Value = 2
$ws.Range("G14").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("H14").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = -19.354838709677
$ws.Range("L16").Value = 78.571428571428
$ws.Range("M16").Value = -41.860465116279
$ws.Range("N16").Value = -85.632183908046
$ws.Range("C17").Value = 7
$ws.Range("G14").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 1
$ws.Range("H14").Copy($ws.Range("E17"))
$ws.Range("E17").Value = 600
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 180
$ws.Range("I17").Value = 81
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = 15.714285714285
$ws.Range("L17").Value = 5.194805194805
$ws.Range("M17").Value = 37.28813559322
$ws.Range("N17").Value = -40.875912408759
$ws.Range("C18").Value = 2
$ws.Range("G14").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1
$ws.Range("H14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -12.5
$ws.Range("I18").Value = 29
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = -25.641025641025
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -25.641025641025
$ws.Range("N18").Value = -90.068493150684
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -66.666666666666
$ws.Range("I19").Value = 120
$ws.Range("J19").Value = 104
$ws.Range("K19").Value = 15.384615384615
$ws.Range("L19").Value = 39.53488372093
$ws.Range("M19").Value = 57.894736842105
$ws.Range("N19").Value = -9.774436090225
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = -42.222222222222
$ws.Range("M20").Value = 44.444444444444
$ws.Range("N20").Value = -88.841201716738
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 62.5
$ws.Range("F21").Value = 36
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = -21.739130434782
$ws.Range("I21").Value = 283
$ws.Range("J21").Value = 298
$ws.Range("K21").Value = -5.033557046979
$ws.Range("L21").Value = 25.777777777777
$ws.Range("M21").Value = 18.410041841004
$ws.Range("N21").Value = -71.210579857578
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("M23").Value = -55.555555555555
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 15
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 316
$ws.Range("J24").Value = 362
$ws.Range("K24").Value = -12.707182320442
$ws.Range("L24").Value = 22.480620155038
$ws.Range("M24").Value = 59.595959595959
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -57.142857142857
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -42.105263157894
$ws.Range("I25").Value = 153
$ws.Range("J25").Value = 133
$ws.Range("K25").Value = 15.037593984962
$ws.Range("L25").Value = 25.409836065573
$ws.Range("M25").Value = -35.443037974683
$ws.Range("G14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("G14").Copy($ws.Range("F26"))
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 5
$ws.Range("K26").Value = -44.444444444444
$ws.Range("L26").Value = -50
$ws.Range("G14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("G14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("H14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("G14").Copy($ws.Range("G27"))
$ws.Range("G27").Value = 1
$ws.Range("H14").Copy($ws.Range("H27"))
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = 28.571428571428
$ws.Range("G28").Value = 2
$ws.Range("G29").Value = 1
